# ARKCORR-18 Added business process definitions for the on enter queue rule.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the rule-table "template" row (row 17): condition / action expressions ---
$ws.Range("C17").Value = '((CaseFilePipelineContext) $model.getPipelineContext()).getEnqueueName().equals("$param")'
$ws.Range("D17").Value = '$model.setBusinessProcessName("$param");'

# --- Update the default-rule row (row 18): rename rule & clarify description ---
$ws.Range("C18").Value = "Entering Queue Name"
$ws.Range("D18").Value = "The business process to be executed when entering that queue"

# --- Apply the existing "data row" look (wrap text, thin border, unlocked) used by B18 ---
# to the new rule rows B19:D23 before filling them in, matching the rest of the rule table.
$ws.Range("B18").Copy()
$ws.Range("B19:D23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Populate the five per-queue business process rules ---
$ws.Range("B19").Value = "Intake queue"
$ws.Range("C19").Value = "Intake"
$ws.Range("D19").Value = "correspondence-extension-intake-process"

$ws.Range("B20").Value = "Fulfill queue"
$ws.Range("C20").Value = "Fulfill"
$ws.Range("D20").Value = "correspondence-extension-fulfill-process"

$ws.Range("B21").Value = "Supervisor Approval queue"
$ws.Range("C21").Value = "Supervisor Approval"
$ws.Range("D21").Value = "correspondence-extension-supervisor-approval-process"

$ws.Range("B22").Value = "Executive Approval queue"
$ws.Range("C22").Value = "Executive Approval"
$ws.Range("D22").Value = "correspondence-extension-executive-approval-process"

$ws.Range("B23").Value = "Release queue"
$ws.Range("C23").Value = "Release"
$ws.Range("D23").Value = "correspondence-extension-release-process"

# --- Move the selection to the last-edited cell, like the author would have left it ---
$ws.Activate()
$ws.Range("D23").Select()
